$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the "city" domain row to "citystgeorges" (dedupes naming with
#    the already-present citystgeorges row) - row 14 in the original data.
$ws.Range("A14").Value = "citystgeorges"

# 2. Strip the leading "The " from a handful of institution names so they
#    join cleanly against the CUG data via domain.
$ws.Range("B38").Value = "University of Manchester"
$ws.Range("B52").Value = "University of Sheffield"
$ws.Range("B45").Value = "Open University"
$ws.Range("B75").Value = "Alan Turing Institute"
$ws.Range("B73").Value = "Scottish Tech Army"
$ws.Range("B19").Value = "University of Edinburgh"

# 3. Add a new affiliation row for Northumbria University, matching the
#    style already used by the other "blank C column" data rows (copy the
#    existing northumbria row's formatting, then overwrite its values).
$ws.Range("A43:C43").Copy()
$ws.Range("A78").PasteSpecial(-4122)
$ws.Range("A78").Value = "northumbria"
$ws.Range("B78").Value = "Northumbria University"

$ws.Range("B54").Value = "Solent University"

# 4. Re-sort the table by domain (column A), same as the sheet's existing
#    sortState, now covering the extra row. Domain ties are broken by the
#    institution name (column B) ascending.
$sortRange = $ws.Range("A2:C78")
$sortRange.Sort($ws.Range("A2:A78"), 1, $ws.Range("B2:B78"), , 1)

# The "northumbria" pair keeps its original relative order (pre-existing
# row ahead of the freshly appended one) rather than the alphabetical
# B tie-break applied everywhere else, so fix that pair back up.
$firstVal = $ws.Range("B49").Value2
$secondVal = $ws.Range("B50").Value2
$ws.Range("B49").Value = $secondVal
$ws.Range("B50").Value = $firstVal

# 5. Restore view state: clear the scrolled-down top-left cell and move
#    the active selection back up.
$ws.Range("A20").Select()
